# Bio data sample.xlsx edits
# - Rename header "Department Code" (E1) to "Department"
# - Fill in D12 with "Ogochukwu" (Other Names for the last row)
# - Reset selection to A1 (no longer C12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Department"
$ws.Range("D12").Value = "Ogochukwu"

$ws.Range("A1").Select()
